$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 6.314527000000001
$ws.Cells.Item(2, 8).Value = 18.943581
$ws.Cells.Item(2, 9).Value = 0.2616724966426195
$ws.Cells.Item(2, 10).Value = 0.2616724966426195
$ws.Cells.Item(2, 13).Value = 6.229715
$ws.Cells.Item(2, 14).Value = 18.689145
$ws.Cells.Item(2, 15).Value = 0.2690007335750987
$ws.Cells.Item(2, 16).Value = 0.2690007335750986
$ws.Cells.Item(2, 17).Value = 39.337703569805
$ws.Cells.Item(2, 18).Value = 354.039332128245
$ws.Cells.Item(2, 19).Value = 0.0703900935532922
$ws.Cells.Item(2, 20).Value = 0.07039009355329219

# Row 3
$ws.Cells.Item(3, 7).Value = 6.314527000000001
$ws.Cells.Item(3, 8).Value = 18.943581
$ws.Cells.Item(3, 9).Value = 0.2616724966426195
$ws.Cells.Item(3, 10).Value = 0.2616724966426195
$ws.Cells.Item(3, 15).Value = 0.5269588201056402
$ws.Cells.Item(3, 16).Value = 0.5269588201056401
$ws.Cells.Item(3, 17).Value = 77.06056999663434
$ws.Cells.Item(3, 18).Value = 693.545129969709
$ws.Cells.Item(3, 19).Value = 0.1378906300848919
$ws.Cells.Item(3, 20).Value = 0.1378906300848919

# Row 4
$ws.Cells.Item(4, 7).Value = 6.314527000000001
$ws.Cells.Item(4, 8).Value = 18.943581
$ws.Cells.Item(4, 9).Value = 0.2616724966426195
$ws.Cells.Item(4, 10).Value = 0.2616724966426195
$ws.Cells.Item(4, 13).Value = 2.383077666666667
$ws.Cells.Item(4, 14).Value = 7.149233000000001
$ws.Cells.Item(4, 15).Value = 0.1029019209546132
$ws.Cells.Item(4, 16).Value = 0.1029019209546131
$ws.Cells.Item(4, 17).Value = 15.04800826926367
$ws.Cells.Item(4, 18).Value = 135.432074423373
$ws.Cells.Item(4, 19).Value = 0.02692660256551511
$ws.Cells.Item(4, 20).Value = 0.02692660256551511

# Row 5
$ws.Cells.Item(5, 7).Value = 6.314527000000001
$ws.Cells.Item(5, 8).Value = 18.943581
$ws.Cells.Item(5, 9).Value = 0.2616724966426195
$ws.Cells.Item(5, 10).Value = 0.2616724966426195
$ws.Cells.Item(5, 13).Value = 2.342239666666666
$ws.Cells.Item(5, 14).Value = 7.026719
$ws.Cells.Item(5, 15).Value = 0.101138525364648
$ws.Cells.Item(5, 16).Value = 0.1011385253646479
$ws.Cells.Item(5, 17).Value = 14.79013561563767
$ws.Cells.Item(5, 18).Value = 133.111220540739
$ws.Cells.Item(5, 19).Value = 0.02646517043892034
$ws.Cells.Item(5, 20).Value = 0.02646517043892033

# Row 6
$ws.Cells.Item(6, 9).Value = 0.1461016137776048
$ws.Cells.Item(6, 10).Value = 0.1461016137776048
$ws.Cells.Item(6, 13).Value = 6.229715
$ws.Cells.Item(6, 14).Value = 18.689145
$ws.Cells.Item(6, 15).Value = 0.2690007335750987
$ws.Cells.Item(6, 16).Value = 0.2690007335750986
$ws.Cells.Item(6, 17).Value = 21.96372200974167
$ws.Cells.Item(6, 18).Value = 197.673498087675
$ws.Cells.Item(6, 19).Value = 0.03930144128268143
$ws.Cells.Item(6, 20).Value = 0.03930144128268142

# Row 7
$ws.Cells.Item(7, 9).Value = 0.1461016137776048
$ws.Cells.Item(7, 10).Value = 0.1461016137776048
$ws.Cells.Item(7, 15).Value = 0.5269588201056402
$ws.Cells.Item(7, 16).Value = 0.5269588201056401
$ws.Cells.Item(7, 19).Value = 0.07698953401177656
$ws.Cells.Item(7, 20).Value = 0.07698953401177655

# Row 8
$ws.Cells.Item(8, 9).Value = 0.1461016137776048
$ws.Cells.Item(8, 10).Value = 0.1461016137776048
$ws.Cells.Item(8, 13).Value = 2.383077666666667
$ws.Cells.Item(8, 14).Value = 7.149233000000001
$ws.Cells.Item(8, 15).Value = 0.1029019209546132
$ws.Cells.Item(8, 16).Value = 0.1029019209546131
$ws.Cells.Item(8, 17).Value = 8.401869972910555
$ws.Cells.Item(8, 18).Value = 75.616829756195
$ws.Cells.Item(8, 19).Value = 0.01503413671228451
$ws.Cells.Item(8, 20).Value = 0.0150341367122845

# Row 9
$ws.Cells.Item(9, 9).Value = 0.1461016137776048
$ws.Cells.Item(9, 10).Value = 0.1461016137776048
$ws.Cells.Item(9, 13).Value = 2.342239666666666
$ws.Cells.Item(9, 14).Value = 7.026719
$ws.Cells.Item(9, 15).Value = 0.101138525364648
$ws.Cells.Item(9, 16).Value = 0.1011385253646479
$ws.Cells.Item(9, 17).Value = 8.257889954653889
$ws.Cells.Item(9, 18).Value = 74.321009591885
$ws.Cells.Item(9, 19).Value = 0.01477650177086228
$ws.Cells.Item(9, 20).Value = 0.01477650177086228

# Row 10
$ws.Cells.Item(10, 7).Value = 0.510814
$ws.Cells.Item(10, 8).Value = 1.532442
$ws.Cells.Item(10, 9).Value = 0.02116801063642661
$ws.Cells.Item(10, 10).Value = 0.02116801063642662
$ws.Cells.Item(10, 13).Value = 6.229715
$ws.Cells.Item(10, 14).Value = 18.689145
$ws.Cells.Item(10, 15).Value = 0.2690007335750987
$ws.Cells.Item(10, 16).Value = 0.2690007335750986
$ws.Cells.Item(10, 17).Value = 3.18222563801
$ws.Cells.Item(10, 18).Value = 28.64003074209
$ws.Cells.Item(10, 19).Value = 0.00569421038952425
$ws.Cells.Item(10, 20).Value = 0.00569421038952425

# Row 11
$ws.Cells.Item(11, 7).Value = 0.510814
$ws.Cells.Item(11, 8).Value = 1.532442
$ws.Cells.Item(11, 9).Value = 0.02116801063642661
$ws.Cells.Item(11, 10).Value = 0.02116801063642662
$ws.Cells.Item(11, 15).Value = 0.5269588201056402
$ws.Cells.Item(11, 16).Value = 0.5269588201056401
$ws.Cells.Item(11, 17).Value = 6.233818938815333
$ws.Cells.Item(11, 18).Value = 56.104370449338
$ws.Cells.Item(11, 19).Value = 0.01115466990895501
$ws.Cells.Item(11, 20).Value = 0.01115466990895501

# Row 12
$ws.Cells.Item(12, 7).Value = 0.510814
$ws.Cells.Item(12, 8).Value = 1.532442
$ws.Cells.Item(12, 9).Value = 0.02116801063642661
$ws.Cells.Item(12, 10).Value = 0.02116801063642662
$ws.Cells.Item(12, 13).Value = 2.383077666666667
$ws.Cells.Item(12, 14).Value = 7.149233000000001
$ws.Cells.Item(12, 15).Value = 0.1029019209546132
$ws.Cells.Item(12, 16).Value = 0.1029019209546131
$ws.Cells.Item(12, 17).Value = 1.217309435220667
$ws.Cells.Item(12, 18).Value = 10.955784916986
$ws.Cells.Item(12, 19).Value = 0.002178228957275982
$ws.Cells.Item(12, 20).Value = 0.002178228957275982

# Row 13
$ws.Cells.Item(13, 7).Value = 0.510814
$ws.Cells.Item(13, 8).Value = 1.532442
$ws.Cells.Item(13, 9).Value = 0.02116801063642661
$ws.Cells.Item(13, 10).Value = 0.02116801063642662
$ws.Cells.Item(13, 13).Value = 2.342239666666666
$ws.Cells.Item(13, 14).Value = 7.026719
$ws.Cells.Item(13, 15).Value = 0.101138525364648
$ws.Cells.Item(13, 16).Value = 0.1011385253646479
$ws.Cells.Item(13, 17).Value = 1.196448813088667
$ws.Cells.Item(13, 18).Value = 10.768039317798
$ws.Cells.Item(13, 19).Value = 0.002140901380671371
$ws.Cells.Item(13, 20).Value = 0.002140901380671371

# Row 14
$ws.Cells.Item(14, 7).Value = 13.78043333333333
$ws.Cells.Item(14, 8).Value = 41.3413
$ws.Cells.Item(14, 9).Value = 0.571057878943349
$ws.Cells.Item(14, 10).Value = 0.5710578789433491
$ws.Cells.Item(14, 13).Value = 6.229715
$ws.Cells.Item(14, 14).Value = 18.689145
$ws.Cells.Item(14, 15).Value = 0.2690007335750987
$ws.Cells.Item(14, 16).Value = 0.2690007335750986
$ws.Cells.Item(14, 17).Value = 85.84817224316667
$ws.Cells.Item(14, 18).Value = 772.6335501885001
$ws.Cells.Item(14, 19).Value = 0.1536149883496008
$ws.Cells.Item(14, 20).Value = 0.1536149883496008

# Row 15
$ws.Cells.Item(15, 7).Value = 13.78043333333333
$ws.Cells.Item(15, 8).Value = 41.3413
$ws.Cells.Item(15, 9).Value = 0.571057878943349
$ws.Cells.Item(15, 10).Value = 0.5710578789433491
$ws.Cells.Item(15, 15).Value = 0.5269588201056402
$ws.Cells.Item(15, 16).Value = 0.5269588201056401
$ws.Cells.Item(15, 17).Value = 168.1722237417445
$ws.Cells.Item(15, 18).Value = 1513.5500136757
$ws.Cells.Item(15, 19).Value = 0.3009239861000167
$ws.Cells.Item(15, 20).Value = 0.3009239861000167

# Row 16
$ws.Cells.Item(16, 7).Value = 13.78043333333333
$ws.Cells.Item(16, 8).Value = 41.3413
$ws.Cells.Item(16, 9).Value = 0.571057878943349
$ws.Cells.Item(16, 10).Value = 0.5710578789433491
$ws.Cells.Item(16, 13).Value = 2.383077666666667
$ws.Cells.Item(16, 14).Value = 7.149233000000001
$ws.Cells.Item(16, 15).Value = 0.1029019209546132
$ws.Cells.Item(16, 16).Value = 0.1029019209546131
$ws.Cells.Item(16, 17).Value = 32.83984291365556
$ws.Cells.Item(16, 18).Value = 295.5585862229
$ws.Cells.Item(16, 19).Value = 0.05876295271953755
$ws.Cells.Item(16, 20).Value = 0.05876295271953755

# Row 17
$ws.Cells.Item(17, 7).Value = 13.78043333333333
$ws.Cells.Item(17, 8).Value = 41.3413
$ws.Cells.Item(17, 9).Value = 0.571057878943349
$ws.Cells.Item(17, 10).Value = 0.5710578789433491
$ws.Cells.Item(17, 13).Value = 2.342239666666666
$ws.Cells.Item(17, 14).Value = 7.026719
$ws.Cells.Item(17, 15).Value = 0.101138525364648
$ws.Cells.Item(17, 16).Value = 0.1011385253646479
$ws.Cells.Item(17, 17).Value = 32.27707757718889
$ws.Cells.Item(17, 18).Value = 290.4936981947
$ws.Cells.Item(17, 19).Value = 0.05775595177419398
$ws.Cells.Item(17, 20).Value = 0.05775595177419397

Write-Output "Done updating Il15-Il15ra sheet values"
